$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the "Category of Evidence" header (column I, row 1) to
# "Study Type/GVD Chapter".
$ws.Range("I1").Value = "Study Type/GVD Chapter"

# Leave the cursor on the edited cell, matching the saved selection state.
$ws.Range("I1").Select() | Out-Null
